# Applies:
#   1) The table on slide 5 switches from the custom "Table_0" style
#      ({A5D22077-A8EC-44D5-8CF5-6E0F30D08CD3}) to the built-in
#      "No Style, Table Grid" style ({777C92A0-F167-4AC2-801E-C76878F47869}).
#   2) The presentation's active colour theme flips from the "Integral"
#      (Red Violet) palette back to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{777C92A0-F167-4AC2-801E-C76878F47869}")
    }
}

# --- 2) Theme colours: Integral/Red Violet -> Office ----------------------
# RGB() isn't available in this host, so the colours below are expressed
# as the packed 0x00BBGGRR integer PowerPoint stores on ColorFormat.RGB.
#   index : theme slot   : office hex : packed value
#     1   : dk1           000000        0
#     2   : lt1           FFFFFF        16777215
#     3   : dk2           44546A        6968388
#     4   : lt2           E7E6E6        15132391
#     5   : accent1       5B9BD5        13998939
#     6   : accent2       ED7D31        3243501
#     7   : accent3       A5A5A5        10855845
#     8   : accent4       FFC000        49407
#     9   : accent5       4472C4        12874308
#    10   : accent6       70AD47        4697456
#    11   : hlink         0563C1        12673797
#    12   : folHlink      954F72        7491477
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$tcs = $slide.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
